# Updated testing of ValidLogin and InValidLogin
# Adds a new "IssueDate" column and inserts a new login-test row
# (xyz/abcd) above the existing admin/manager/Enter row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2, pushing the existing
# admin / manager / Enter row down to row 3.
$ws.Rows("2:2").Insert()

# Row 1: new header cell for the IssueDate column.
$ws.Range("D1").Value = "IssueDate"

# Row 2: new test data (xyz / abcd) plus two date-formatted,
# otherwise empty cells under the new IssueDate column.
$ws.Range("A2").Value = "xyz"
$ws.Range("B2").Value = "abcd"
$ws.Range("D2").NumberFormat = "d-mmm-yy"
$ws.Range("E2").NumberFormat = "d-mmm-yy"

# Row 3 (previously row 2): fill in the "Enter" value that was
# missing from column C for the admin/manager row.
$ws.Range("C3").Value = "Enter"
